$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 248
$ws.Range("D248").Value = 44795
$ws.Range("J248").Value = 2000
$ws.Range("K248").Value = 1300
$ws.Range("L248").Value = 1300
$ws.Range("M248").Value = 1300
$ws.Range("O248").Value = "Región del Maule"
$ws.Range("P248").Value = 1300

# Row 249
$ws.Range("D249").Value = 44795
$ws.Range("I249").Value = "Segunda"

# Row 250
$ws.Range("D250").Value = 44508
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 6000
$ws.Range("K250").Value = 700
$ws.Range("L250").Value = 700
$ws.Range("M250").Value = 700
$ws.Range("O250").Value = "Provincia del Elquí"
$ws.Range("P250").Value = 700

# Row 251
$ws.Range("D251").Value = 44775
$ws.Range("K251").Value = 1000
$ws.Range("L251").Value = 1000
$ws.Range("M251").Value = 1000
$ws.Range("P251").Value = 1000

# Row 252
$ws.Range("D252").Value = 44775
$ws.Range("I252").Value = "Segunda"
$ws.Range("K252").Value = 800
$ws.Range("L252").Value = 800
$ws.Range("M252").Value = 800
$ws.Range("P252").Value = 800

# Row 253
$ws.Range("D253").Value = 44335
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 3000
$ws.Range("K253").Value = 550
$ws.Range("L253").Value = 550
$ws.Range("M253").Value = 550
$ws.Range("P253").Value = 550

# Row 254
$ws.Range("D254").Value = 44448
$ws.Range("J254").Value = 3000
$ws.Range("K254").Value = 500
$ws.Range("L254").Value = 500
$ws.Range("M254").Value = 500
$ws.Range("P254").Value = 500

# Row 255
$ws.Range("D255").Value = 44627
$ws.Range("I255").Value = "Segunda"
$ws.Range("J255").Value = 2000
$ws.Range("K255").Value = 1000
$ws.Range("L255").Value = 1000
$ws.Range("M255").Value = 1000
$ws.Range("P255").Value = 1000

# Row 256
$ws.Range("D256").Value = 44315
$ws.Range("J256").Value = 4000
$ws.Range("K256").Value = 400
$ws.Range("L256").Value = 400
$ws.Range("M256").Value = 400
$ws.Range("P256").Value = 400

# Row 257
$ws.Range("D257").Value = 44386
$ws.Range("J257").Value = 8000
$ws.Range("K257").Value = 500
$ws.Range("L257").Value = 500
$ws.Range("M257").Value = 500
$ws.Range("P257").Value = 500

# Row 258
$ws.Range("D258").Value = 44327
$ws.Range("J258").Value = 5000
$ws.Range("K258").Value = 450
$ws.Range("L258").Value = 450
$ws.Range("M258").Value = 450
$ws.Range("P258").Value = 450

# Row 259
$ws.Range("D259").Value = 44316
$ws.Range("J259").Value = 5000
$ws.Range("K259").Value = 400
$ws.Range("L259").Value = 400
$ws.Range("M259").Value = 400
$ws.Range("O259").Value = "Región del Maule"
$ws.Range("P259").Value = 400

# Row 260
$ws.Range("D260").Value = 44264
$ws.Range("J260").Value = 3000
$ws.Range("K260").Value = 1000
$ws.Range("L260").Value = 1000
$ws.Range("M260").Value = 1000
$ws.Range("P260").Value = 1000

# Row 261
$ws.Range("D261").Value = 44523
$ws.Range("J261").Value = 6000
$ws.Range("K261").Value = 800
$ws.Range("L261").Value = 800
$ws.Range("M261").Value = 800
$ws.Range("O261").Value = "Provincia del Elquí"
$ws.Range("P261").Value = 800

# Row 262
$ws.Range("D262").Value = 44320
$ws.Range("J262").Value = 4000
$ws.Range("K262").Value = 500
$ws.Range("L262").Value = 500
$ws.Range("M262").Value = 500
$ws.Range("P262").Value = 500

# Row 263
$ws.Range("D263").Value = 44566
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 3000
$ws.Range("K263").Value = 700
$ws.Range("L263").Value = 700
$ws.Range("M263").Value = 700
$ws.Range("P263").Value = 700

# Row 264
$ws.Range("D264").Value = 44764
$ws.Range("J264").Value = 5000
$ws.Range("K264").Value = 1300
$ws.Range("L264").Value = 1500
$ws.Range("M264").Value = 1420
$ws.Range("P264").Value = 1420

# Row 265
$ws.Range("D265").Value = 44764
$ws.Range("I265").Value = "Segunda"
$ws.Range("J265").Value = 2000
$ws.Range("K265").Value = 1000
$ws.Range("L265").Value = 1000
$ws.Range("M265").Value = 1000
$ws.Range("P265").Value = 1000

# Row 266
$ws.Range("D266").Value = 44552
$ws.Range("J266").Value = 3000

# Row 267
$ws.Range("D267").Value = 44231
$ws.Range("J267").Value = 3000
$ws.Range("K267").Value = 800
$ws.Range("L267").Value = 800
$ws.Range("M267").Value = 800
$ws.Range("P267").Value = 800

# Row 268
$ws.Range("D268").Value = 44557
$ws.Range("J268").Value = 5000
$ws.Range("K268").Value = 600
$ws.Range("L268").Value = 600
$ws.Range("M268").Value = 600
$ws.Range("O268").Value = "Región del Maule"
$ws.Range("P268").Value = 600

# Row 269
$ws.Range("D269").Value = 44711
$ws.Range("J269").Value = 5000
$ws.Range("K269").Value = 1000
$ws.Range("L269").Value = 1000
$ws.Range("M269").Value = 1000
$ws.Range("P269").Value = 1000

# Row 270
$ws.Range("D270").Value = 44490
$ws.Range("K270").Value = 700
$ws.Range("L270").Value = 700
$ws.Range("M270").Value = 700
$ws.Range("O270").Value = "Provincia del Elquí"
$ws.Range("P270").Value = 700

# Row 271
$ws.Range("D271").Value = 44778
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 2500
$ws.Range("K271").Value = 1200
$ws.Range("L271").Value = 1200
$ws.Range("M271").Value = 1200
$ws.Range("P271").Value = 1200

# Row 272
$ws.Range("D272").Value = 44279
$ws.Range("J272").Value = 3000
$ws.Range("K272").Value = 800
$ws.Range("L272").Value = 800
$ws.Range("M272").Value = 800
$ws.Range("P272").Value = 800

# Row 273
$ws.Range("D273").Value = 44431
$ws.Range("K273").Value = 350
$ws.Range("L273").Value = 350
$ws.Range("M273").Value = 350
$ws.Range("P273").Value = 350

# Row 274
$ws.Range("D274").Value = 44749
$ws.Range("J274").Value = 2000
$ws.Range("K274").Value = 1200
$ws.Range("L274").Value = 1200
$ws.Range("M274").Value = 1200
$ws.Range("P274").Value = 1200

# Row 275
$ws.Range("D275").Value = 44749
$ws.Range("I275").Value = "Segunda"
$ws.Range("J275").Value = 3000
$ws.Range("K275").Value = 900
$ws.Range("L275").Value = 900
$ws.Range("M275").Value = 900
$ws.Range("P275").Value = 900

# Row 276
$ws.Range("D276").Value = 44727
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 3000
$ws.Range("K276").Value = 1100
$ws.Range("L276").Value = 1100
$ws.Range("M276").Value = 1100
$ws.Range("P276").Value = 1100

# Row 277
$ws.Range("D277").Value = 44761
$ws.Range("J277").Value = 2000
$ws.Range("K277").Value = 1500
$ws.Range("L277").Value = 1500
$ws.Range("M277").Value = 1500
$ws.Range("P277").Value = 1500

# Row 278
$ws.Range("D278").Value = 44761
$ws.Range("I278").Value = "Segunda"
$ws.Range("J278").Value = 1200
$ws.Range("K278").Value = 1000
$ws.Range("L278").Value = 1000
$ws.Range("M278").Value = 1000
$ws.Range("P278").Value = 1000

# Row 279
$ws.Range("D279").Value = 44223
$ws.Range("K279").Value = 800
$ws.Range("L279").Value = 800
$ws.Range("M279").Value = 800
$ws.Range("P279").Value = 800

# Row 280
$ws.Range("D280").Value = 44342
$ws.Range("I280").Value = "Primera"
$ws.Range("J280").Value = 5000
$ws.Range("K280").Value = 600
$ws.Range("L280").Value = 600
$ws.Range("M280").Value = 600
$ws.Range("P280").Value = 600

# Row 281
$ws.Range("D281").Value = 44678
$ws.Range("K281").Value = 1000
$ws.Range("L281").Value = 1000
$ws.Range("M281").Value = 1000
$ws.Range("P281").Value = 1000

# Row 282
$ws.Range("D282").Value = 44678
$ws.Range("I282").Value = "Segunda"
$ws.Range("J282").Value = 2000
$ws.Range("K282").Value = 800
$ws.Range("L282").Value = 800
$ws.Range("M282").Value = 800
$ws.Range("P282").Value = 800

# Row 283
$ws.Range("D283").Value = 44551
$ws.Range("I283").Value = "Primera"
$ws.Range("J283").Value = 3000
$ws.Range("K283").Value = 900
$ws.Range("L283").Value = 900
$ws.Range("M283").Value = 900
$ws.Range("P283").Value = 900

# Row 284
$ws.Range("D284").Value = 44291
$ws.Range("K284").Value = 700
$ws.Range("L284").Value = 700
$ws.Range("M284").Value = 700
$ws.Range("P284").Value = 700

# Row 285
$ws.Range("D285").Value = 44629
$ws.Range("K285").Value = 1000
$ws.Range("L285").Value = 1000
$ws.Range("M285").Value = 1000
$ws.Range("P285").Value = 1000

# Row 286
$ws.Range("D286").Value = 44449
$ws.Range("J286").Value = 3000

# Row 287
$ws.Range("D287").Value = 44449
$ws.Range("I287").Value = "Segunda"
$ws.Range("J287").Value = 2000
$ws.Range("K287").Value = 300
$ws.Range("L287").Value = 300
$ws.Range("M287").Value = 300
$ws.Range("P287").Value = 300

# Row 288
$ws.Range("D288").Value = 44376
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 4000
$ws.Range("K288").Value = 500
$ws.Range("L288").Value = 500
$ws.Range("M288").Value = 500
$ws.Range("P288").Value = 500

# Row 289
$ws.Range("D289").Value = 44391
$ws.Range("K289").Value = 400
$ws.Range("L289").Value = 400
$ws.Range("M289").Value = 400
$ws.Range("P289").Value = 400

# Row 290
$ws.Range("D290").Value = 44600
$ws.Range("I290").Value = "Segunda"
$ws.Range("J290").Value = 2000
$ws.Range("K290").Value = 800
$ws.Range("L290").Value = 800
$ws.Range("M290").Value = 800
$ws.Range("P290").Value = 800

# Row 291
$ws.Range("D291").Value = 44371
$ws.Range("J291").Value = 5000
$ws.Range("K291").Value = 450
$ws.Range("L291").Value = 450
$ws.Range("M291").Value = 450
$ws.Range("P291").Value = 450

# Row 292
$ws.Range("D292").Value = 44355
$ws.Range("J292").Value = 5000
$ws.Range("K292").Value = 450
$ws.Range("L292").Value = 450
$ws.Range("M292").Value = 450
$ws.Range("P292").Value = 450

# Row 293
$ws.Range("D293").Value = 44685
$ws.Range("I293").Value = "Primera"
$ws.Range("K293").Value = 1000
$ws.Range("L293").Value = 1000
$ws.Range("M293").Value = 1000
$ws.Range("P293").Value = 1000

# Row 294
$ws.Range("D294").Value = 44434
$ws.Range("J294").Value = 2000

# Row 295
$ws.Range("D295").Value = 44434
$ws.Range("I295").Value = "Segunda"
$ws.Range("K295").Value = 350
$ws.Range("L295").Value = 350
$ws.Range("M295").Value = 350
$ws.Range("P295").Value = 350

# Row 296
$ws.Range("D296").Value = 44462
$ws.Range("J296").Value = 3000
$ws.Range("K296").Value = 500
$ws.Range("L296").Value = 500
$ws.Range("M296").Value = 500
$ws.Range("P296").Value = 500

# Row 297
$ws.Range("D297").Value = 44278
$ws.Range("K297").Value = 800
$ws.Range("L297").Value = 800
$ws.Range("M297").Value = 800
$ws.Range("P297").Value = 800

# Row 298
$ws.Range("D298").Value = 44272
$ws.Range("I298").Value = "Primera"
$ws.Range("K298").Value = 800
$ws.Range("L298").Value = 800
$ws.Range("M298").Value = 800
$ws.Range("P298").Value = 800

# Row 299
$ws.Range("D299").Value = 44781
$ws.Range("J299").Value = 3000
$ws.Range("K299").Value = 1300
$ws.Range("L299").Value = 1300
$ws.Range("M299").Value = 1300
$ws.Range("P299").Value = 1300

# Row 300
$ws.Range("D300").Value = 44781
$ws.Range("K300").Value = 1000
$ws.Range("L300").Value = 1000
$ws.Range("M300").Value = 1000
$ws.Range("P300").Value = 1000

# Row 301
$ws.Range("D301").Value = 44336
$ws.Range("J301").Value = 6000

# Row 302
$ws.Range("D302").Value = 44336

# Row 303
$ws.Range("D303").Value = 44343
$ws.Range("J303").Value = 2000
$ws.Range("K303").Value = 600
$ws.Range("L303").Value = 600
$ws.Range("M303").Value = 600
$ws.Range("P303").Value = 600

# Row 304
$ws.Range("D304").Value = 44343
$ws.Range("I304").Value = "Segunda"
$ws.Range("J304").Value = 2000

# Row 305
$ws.Range("D305").Value = 44365
$ws.Range("I305").Value = "Primera"
$ws.Range("K305").Value = 400
$ws.Range("L305").Value = 400
$ws.Range("M305").Value = 400
$ws.Range("P305").Value = 400

# Row 306
$ws.Range("D306").Value = 44421
$ws.Range("J306").Value = 3000
$ws.Range("K306").Value = 400
$ws.Range("L306").Value = 400
$ws.Range("M306").Value = 400
$ws.Range("P306").Value = 400

# Row 307
$ws.Range("D307").Value = 44421
$ws.Range("J307").Value = 3000
$ws.Range("K307").Value = 250
$ws.Range("L307").Value = 250
$ws.Range("M307").Value = 250
$ws.Range("P307").Value = 250

# Row 308
$ws.Range("D308").Value = 44419
$ws.Range("K308").Value = 500
$ws.Range("L308").Value = 500
$ws.Range("M308").Value = 500
$ws.Range("P308").Value = 500

# Row 309
$ws.Range("D309").Value = 44419
$ws.Range("I309").Value = "Segunda"
$ws.Range("J309").Value = 4000
$ws.Range("K309").Value = 350
$ws.Range("L309").Value = 350
$ws.Range("M309").Value = 350
$ws.Range("P309").Value = 350

# Row 310
$ws.Range("D310").Value = 44669
$ws.Range("I310").Value = "Primera"
$ws.Range("K310").Value = 1000
$ws.Range("L310").Value = 1000
$ws.Range("M310").Value = 1000
$ws.Range("P310").Value = 1000

# Row 311
$ws.Range("D311").Value = 44699
$ws.Range("J311").Value = 3000
$ws.Range("K311").Value = 1000
$ws.Range("L311").Value = 1000
$ws.Range("M311").Value = 1000
$ws.Range("P311").Value = 1000

# Row 312
$ws.Range("D312").Value = 44699
$ws.Range("J312").Value = 2000
$ws.Range("K312").Value = 800
$ws.Range("L312").Value = 800
$ws.Range("M312").Value = 800
$ws.Range("P312").Value = 800

# Row 313
$ws.Range("D313").Value = 44405
$ws.Range("J313").Value = 2000
$ws.Range("K313").Value = 500
$ws.Range("L313").Value = 500
$ws.Range("M313").Value = 500
$ws.Range("P313").Value = 500

# Row 314
$ws.Range("D314").Value = 44405
$ws.Range("I314").Value = "Segunda"
$ws.Range("J314").Value = 4000
$ws.Range("K314").Value = 350
$ws.Range("L314").Value = 350
$ws.Range("M314").Value = 350
$ws.Range("P314").Value = 350

# Row 315
$ws.Range("D315").Value = 44273
$ws.Range("I315").Value = "Primera"

# Row 316
$ws.Range("D316").Value = 44777
$ws.Range("K316").Value = 1200
$ws.Range("L316").Value = 1200
$ws.Range("M316").Value = 1200
$ws.Range("P316").Value = 1200

# Row 317
$ws.Range("D317").Value = 44777
$ws.Range("I317").Value = "Segunda"
$ws.Range("K317").Value = 800
$ws.Range("L317").Value = 800
$ws.Range("M317").Value = 800
$ws.Range("P317").Value = 800

# Row 318
$ws.Range("D318").Value = 44309
$ws.Range("I318").Value = "Primera"
$ws.Range("K318").Value = 600
$ws.Range("L318").Value = 600
$ws.Range("M318").Value = 600
$ws.Range("P318").Value = 600

# Row 319
$ws.Range("D319").Value = 44771
$ws.Range("K319").Value = 1000
$ws.Range("L319").Value = 1000
$ws.Range("M319").Value = 1000
$ws.Range("P319").Value = 1000

# Row 320
$ws.Range("D320").Value = 44771
$ws.Range("I320").Value = "Segunda"
$ws.Range("K320").Value = 800
$ws.Range("L320").Value = 800
$ws.Range("M320").Value = 800
$ws.Range("P320").Value = 800

# Row 321
$ws.Range("D321").Value = 44267
$ws.Range("I321").Value = "Primera"
$ws.Range("K321").Value = 800
$ws.Range("L321").Value = 800
$ws.Range("M321").Value = 800
$ws.Range("P321").Value = 800

# Row 322
$ws.Range("D322").Value = 44413

# Row 323
$ws.Range("D323").Value = 44413
$ws.Range("J323").Value = 3000

# Row 324
$ws.Range("D324").Value = 44328
$ws.Range("J324").Value = 3000
$ws.Range("K324").Value = 500
$ws.Range("L324").Value = 500
$ws.Range("M324").Value = 500
$ws.Range("P324").Value = 500

# Row 325
$ws.Range("D325").Value = 44328
$ws.Range("H325").Value = "Crespo record"
$ws.Range("I325").Value = "Segunda"
$ws.Range("K325").Value = 350
$ws.Range("L325").Value = 350
$ws.Range("M325").Value = 350
$ws.Range("P325").Value = 350

# Row 326
$ws.Range("D326").Value = 44515
$ws.Range("J326").Value = 4000
$ws.Range("K326").Value = 900
$ws.Range("L326").Value = 900
$ws.Range("M326").Value = 900
$ws.Range("P326").Value = 900

# Row 327
$ws.Range("D327").Value = 44515
$ws.Range("H327").Value = "Morada(o)"
$ws.Range("J327").Value = 2000
$ws.Range("K327").Value = 1000
$ws.Range("L327").Value = 1000
$ws.Range("M327").Value = 1000
$ws.Range("P327").Value = 1000

# Row 328
$ws.Range("D328").Value = 44356
$ws.Range("J328").Value = 5000

# Row 329
$ws.Range("D329").Value = 44379
$ws.Range("I329").Value = "Primera"
$ws.Range("J329").Value = 4000
$ws.Range("K329").Value = 500
$ws.Range("L329").Value = 500
$ws.Range("M329").Value = 500
$ws.Range("P329").Value = 500

# Row 330
$ws.Range("D330").Value = 44322
$ws.Range("J330").Value = 3000
$ws.Range("K330").Value = 450
$ws.Range("L330").Value = 450
$ws.Range("M330").Value = 450
$ws.Range("O330").Value = "Región del Maule"
$ws.Range("P330").Value = 450

# Row 331
$ws.Range("D331").Value = 44322
$ws.Range("I331").Value = "Segunda"
$ws.Range("J331").Value = 3000
$ws.Range("K331").Value = 350
$ws.Range("L331").Value = 350
$ws.Range("M331").Value = 350
$ws.Range("P331").Value = 350

# Row 332
$ws.Range("D332").Value = 44497
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 5000
$ws.Range("K332").Value = 600
$ws.Range("L332").Value = 600
$ws.Range("M332").Value = 600
$ws.Range("O332").Value = "Provincia del Elquí"
$ws.Range("P332").Value = 600

# Row 333
$ws.Range("D333").Value = 44782
$ws.Range("J333").Value = 2500
$ws.Range("K333").Value = 1300
$ws.Range("L333").Value = 1300
$ws.Range("M333").Value = 1300
$ws.Range("P333").Value = 1300

# Row 334
$ws.Range("D334").Value = 44782
$ws.Range("J334").Value = 3000
$ws.Range("K334").Value = 1000
$ws.Range("L334").Value = 1000
$ws.Range("M334").Value = 1000
$ws.Range("P334").Value = 1000

# Row 335
$ws.Range("D335").Value = 44435
$ws.Range("J335").Value = 6000
$ws.Range("K335").Value = 500
$ws.Range("L335").Value = 500
$ws.Range("M335").Value = 500
$ws.Range("P335").Value = 500

# Row 336
$ws.Range("D336").Value = 44435
$ws.Range("I336").Value = "Segunda"
$ws.Range("J336").Value = 15000
$ws.Range("K336").Value = 350
$ws.Range("L336").Value = 350
$ws.Range("M336").Value = 350
$ws.Range("P336").Value = 350

# Row 337
$ws.Range("D337").Value = 44319
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 4000

# Row 338
$ws.Range("D338").Value = 44344
$ws.Range("J338").Value = 3000
$ws.Range("K338").Value = 600
$ws.Range("L338").Value = 600
$ws.Range("M338").Value = 600
$ws.Range("P338").Value = 600

# Row 339
$ws.Range("D339").Value = 44344
$ws.Range("K339").Value = 400
$ws.Range("L339").Value = 400
$ws.Range("M339").Value = 400
$ws.Range("P339").Value = 400

# Row 340
$ws.Range("D340").Value = 44455
$ws.Range("J340").Value = 2000
$ws.Range("K340").Value = 500
$ws.Range("L340").Value = 500
$ws.Range("M340").Value = 500
$ws.Range("O340").Value = "Región del Maule"
$ws.Range("P340").Value = 500

# Row 341
$ws.Range("D341").Value = 44455
$ws.Range("I341").Value = "Segunda"
$ws.Range("J341").Value = 2000
$ws.Range("K341").Value = 300
$ws.Range("L341").Value = 300
$ws.Range("M341").Value = 300
$ws.Range("P341").Value = 300

# Row 342
$ws.Range("D342").Value = 44504
$ws.Range("I342").Value = "Primera"
$ws.Range("J342").Value = 6000
$ws.Range("K342").Value = 600
$ws.Range("L342").Value = 600
$ws.Range("M342").Value = 600
$ws.Range("O342").Value = "Provincia del Elquí"
$ws.Range("P342").Value = 600

# Row 343
$ws.Range("D343").Value = 44484
$ws.Range("J343").Value = 3000
$ws.Range("K343").Value = 700
$ws.Range("L343").Value = 700
$ws.Range("M343").Value = 700
$ws.Range("P343").Value = 700

# Row 344
$ws.Range("D344").Value = 44665
$ws.Range("J344").Value = 2000
$ws.Range("K344").Value = 1000
$ws.Range("L344").Value = 1000
$ws.Range("M344").Value = 1000
$ws.Range("P344").Value = 1000

# Row 345
$ws.Range("D345").Value = 44452
$ws.Range("J345").Value = 2000
$ws.Range("K345").Value = 500
$ws.Range("L345").Value = 500
$ws.Range("M345").Value = 500
$ws.Range("O345").Value = "Región del Maule"
$ws.Range("P345").Value = 500

# Row 346
$ws.Range("D346").Value = 44452
$ws.Range("J346").Value = 5000
$ws.Range("K346").Value = 300
$ws.Range("L346").Value = 300
$ws.Range("M346").Value = 300
$ws.Range("P346").Value = 300

# Row 347
$ws.Range("D347").Value = 44510
$ws.Range("J347").Value = 6000
$ws.Range("K347").Value = 800
$ws.Range("L347").Value = 800
$ws.Range("M347").Value = 800
$ws.Range("O347").Value = "Provincia del Elquí"
$ws.Range("P347").Value = 800

# Row 348
$ws.Range("D348").Value = 44189
$ws.Range("I348").Value = "Segunda"
$ws.Range("K348").Value = 500
$ws.Range("L348").Value = 500
$ws.Range("M348").Value = 500
$ws.Range("P348").Value = 500

# Row 349 (new row)
$ws.Range("A349").Value = 5
$ws.Range("B349").Value = "Macroferia Regional de Talca"
$ws.Range("C349").Value = "Maule"
$ws.Range("D349").Value = 44701
$ws.Range("E349").Value = 7
$ws.Range("F349").Value = 100112006
$ws.Range("G349").Value = "Repollo"
$ws.Range("H349").Value = "Crespo record"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 3000
$ws.Range("K349").Value = 1000
$ws.Range("L349").Value = 1000
$ws.Range("M349").Value = 1000
$ws.Range("N349").Value = "$/unidad"
$ws.Range("O349").Value = "Región del Maule"
$ws.Range("P349").Value = 1000
$ws.Range("Q349").Value = 1
$ws.Range("R349").Value = "Hortaliza"
$ws.Range("D349").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 350 (new row)
$ws.Range("A350").Value = 5
$ws.Range("B350").Value = "Macroferia Regional de Talca"
$ws.Range("C350").Value = "Maule"
$ws.Range("D350").Value = 44516
$ws.Range("E350").Value = 7
$ws.Range("F350").Value = 100112006
$ws.Range("G350").Value = "Repollo"
$ws.Range("H350").Value = "Crespo record"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 3000
$ws.Range("K350").Value = 900
$ws.Range("L350").Value = 900
$ws.Range("M350").Value = 900
$ws.Range("N350").Value = "$/unidad"
$ws.Range("O350").Value = "Región del Maule"
$ws.Range("P350").Value = 900
$ws.Range("Q350").Value = 1
$ws.Range("R350").Value = "Hortaliza"
$ws.Range("D350").NumberFormat = "YYYY-MM-DD HH:MM:SS"
